$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.680.71"
$ws.Range("E2").Value = "  +0.67%  "

$ws.Range("D3").Value = "3.240.28"
$ws.Range("E3").Value = "  +1.47%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").Value = "'605.68"
$ws.Range("E5").Value = "  +1.40%  "

$ws.Range("D6").Value = "'157.39"
$ws.Range("E6").Value = "  +2.43%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "3.238.35"
$ws.Range("E8").Value = "  +1.48%  "

$ws.Range("D9").Value = "'0.547"
$ws.Range("E9").Value = "  +2.20%  "

$ws.Range("D10").Value = "'0.161"
$ws.Range("E10").Value = "  +0.44%  "

$ws.Range("D11").Value = "'5.67"
$ws.Range("E11").Value = "  -7.24%  "

$ws.Range("D12").Value = "'0.511"
$ws.Range("E12").Value = "  -0.69%  "

$ws.Range("D13").Value = "'0.0000271"
$ws.Range("E13").Value = "  +0.96%  "

$ws.Range("D14").Value = "'39.01"
$ws.Range("E14").Value = "  +0.26%  "

$ws.Range("D15").Value = "3.774.32"
$ws.Range("E15").Value = "  +1.52%  "

$ws.Range("D16").Value = "66.694.90"
$ws.Range("E16").Value = "  +0.76%  "

$ws.Range("D17").Value = "'7.51"
$ws.Range("E17").Value = "  +1.04%  "

$ws.Range("D18").Value = "3.239.75"
$ws.Range("E18").Value = "  +1.35%  "

$ws.Range("E19").Value = "  +1.15%  "

$ws.Range("D20").Value = "'512.75"
$ws.Range("E20").Value = "  +0.50%  "

$ws.Range("E21").Value = "  +0.38%  "

$ws.Range("D22").Value = "'0.737"
$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("D23").Value = "'8.11"
$ws.Range("E23").Value = "  +1.42%  "

$ws.Range("E24").Value = "  -1.00%  "

$ws.Range("D25").Value = "'84.85"
$ws.Range("E25").Value = "  -0.07%  "

$ws.Range("E26").Value = "  +0.19%  "

$ws.Range("D27").Value = "'9.33"
$ws.Range("E27").Value = "  +0.70%  "

$ws.Range("E28").Value = "  +0.66%  "

$ws.Range("D29").Value = "'2.42"
$ws.Range("E29").Value = "  +5.69%  "

$ws.Range("D30").Value = "'3.02"
$ws.Range("E30").Value = "  +4.34%  "

$ws.Range("D31").Value = "'7.09"
$ws.Range("E31").Value = "  +1.66%  "

$ws.Range("D32").Value = "'28.32"
$ws.Range("E32").Value = "  +0.28%  "

$ws.Range("E33").Value = "  +0.04%  "

$ws.Range("E34").Value = "  -3.60%  "

$ws.Range("D35").Value = "'6.55"
$ws.Range("E35").Value = "  +0.45%  "

$ws.Range("D36").Value = "'519.48"
$ws.Range("E36").Value = "  +7.26%  "

$ws.Range("D37").Value = "'0.0948"
$ws.Range("E37").Value = "  +5.38%  "

$ws.Range("D38").Value = "'56.14"
$ws.Range("E38").Value = "  +2.40%  "

$ws.Range("D39").Value = "0.0₃0763"
$ws.Range("E39").Value = "  +17.93%  "

$ws.Range("D40").Value = "'0.0421"
$ws.Range("E40").Value = "  +0.84%  "

$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "'3.01"
$ws.Range("E41").Value = "  +3.65%  "

$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.128"
$ws.Range("E42").Value = "  +5.03%  "

$ws.Range("D43").Value = "'8.84"
$ws.Range("E43").Value = "  +0.14%  "

$ws.Range("E44").Value = "  +1.29%  "

$ws.Range("D45").Value = "'2.49"
$ws.Range("E45").Value = "  +3.27%  "

$ws.Range("D46").Value = "2.862.33"
$ws.Range("E46").Value = "  -1.95%  "

$ws.Range("D47").Value = "'28.51"
$ws.Range("E47").Value = "  +0.30%  "

$ws.Range("E48").Value = "  +4.45%  "

$ws.Range("E49").Value = "  -0.04%  "

$ws.Range("E50").Value = "  +0.48%  "

$ws.Range("D51").Value = "'2.63"
$ws.Range("E51").Value = "  +2.07%  "
